$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.950.50'
$ws.Range('E2').Value = '  +2.01%  '
$ws.Range('D3').Value = '1.845.39'
$ws.Range('E3').Value = '  +2.22%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''310.62'
$ws.Range('E5').Value = '  +1.23%  '
$ws.Range('D6').Value = '''1.008'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '''0.4669'
$ws.Range('E7').Value = '  +3.38%  '
$ws.Range('D8').Value = '''0.3629'
$ws.Range('E8').Value = '  +0.96%  '
$ws.Range('D9').Value = '''0.07174'
$ws.Range('E9').Value = '  +1.63%  '
$ws.Range('D10').Value = '''0.9251'
$ws.Range('E10').Value = '  +3.92%  '
$ws.Range('D11').Value = '''19.58'
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('E12').Value = '  -1.51%  '
$ws.Range('D13').Value = '1.855.13'
$ws.Range('E13').Value = '  +2.16%  '
$ws.Range('D14').Value = '''5.299'
$ws.Range('E14').Value = '  +0.49%  '
$ws.Range('D15').Value = '''6.403'
$ws.Range('E15').Value = '  +1.68%  '
$ws.Range('D16').Value = '''88.40'
$ws.Range('E16').Value = '  +3.59%  '
$ws.Range('D17').Value = '''1.011'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = '''0.000008589'
$ws.Range('E18').Value = '  +1.18%  '
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').Value = '26.978.41'
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('D21').Value = '''14.40'
$ws.Range('E21').Value = '  +1.65%  '
$ws.Range('D22').Value = '''5.032'
$ws.Range('E22').Value = '  +1.52%  '
$ws.Range('D23').Value = '''10.63'
$ws.Range('E23').Value = '  +1.24%  '
$ws.Range('D24').Value = '''1.932'
$ws.Range('E24').Value = '  -1.57%  '
$ws.Range('D25').Value = '''152.12'
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  +2.40%  '
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('D28').Value = '''113.95'
$ws.Range('E28').Value = '  +1.68%  '
$ws.Range('D29').Value = '''4.914'
$ws.Range('E29').Value = '  +1.47%  '
$ws.Range('D30').Value = '''0.08858'
$ws.Range('E30').Value = '  +1.97%  '
$ws.Range('D31').Value = '''3.184'
$ws.Range('E31').Value = '  +3.47%  '
$ws.Range('D32').Value = '''2.849'
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('D33').Value = '''1.175'
$ws.Range('E33').Value = '  +6.23%  '
$ws.Range('D34').Value = '''0.7447'
$ws.Range('E34').Value = '  +2.51%  '
$ws.Range('D35').Value = '''4.470'
$ws.Range('E35').Value = '  +0.30%  '
$ws.Range('E36').Value = '  +0.86%  '
$ws.Range('D37').Value = '''2.975'
$ws.Range('E37').Value = '  +2.16%  '
$ws.Range('D38').Value = '''0.01938'
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('D39').Value = '''0.05165'
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('D40').Value = '''0.5143'
$ws.Range('E40').Value = '  +1.72%  '
$ws.Range('D41').Value = '''6.880'
$ws.Range('E41').Value = '  +1.29%  '
$ws.Range('D42').Value = '''0.1510'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').Value = '''8.173'
$ws.Range('E43').Value = '  +1.92%  '
$ws.Range('E44').Value = '  +5.96%  '
$ws.Range('D45').Value = '''0.4704'
$ws.Range('E45').Value = '  +0.80%  '
$ws.Range('D46').Value = '''1.009'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '''100.36'
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('D48').Value = '''1.599'
$ws.Range('E48').Value = '  +1.61%  '
$ws.Range('D49').Value = '''0.06048'
$ws.Range('E49').Value = '  +1.30%  '
$ws.Range('D50').Value = '''64.55'
$ws.Range('E50').Value = '  +1.64%  '
$ws.Range('D51').Value = '''36.10'
$ws.Range('E51').Value = '  +0.03%  '
